# Appium_Page_Object_Model/Excel/appium_data.xlsx
# "Second test case COMPLETED." -> add a VillaTest sheet (city/villa pairs)
# after the existing LoginTest sheet, and leave it as the active tab.

$wb = $excel.ActiveWorkbook

# LoginTest: drop the single-cell selection in favour of the whole used
# range (A1:A3) and stop being the selected/active tab.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:A3").Select() | Out-Null

# Insert the new sheet right after LoginTest.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "VillaTest"

$ws2.Range("A1").Value = "city"
$ws2.Range("B1").Value = "villa"
$ws2.Range("A2").Value = "Delhi"
$ws2.Range("B2").Value = "Hotel Saffron"
$ws2.Range("A3").Value = "Dubai"
$ws2.Range("B3").Value = "Deira Suites"

# Leave the selection on B5, matching the author's last interaction.
$ws2.Range("B5").Select() | Out-Null
